# Scheduled runner update: refresh currentAveragePrice / LevePrice / LeveProfit
# columns (H:N) across several rows in the per-job Leve profit sheets.
# Values below were produced by the external price-refresh job; this script
# simply writes the newly computed figures into the already-open workbook.

$wb = $excel.ActiveWorkbook

function Set-Row {
    param(
        [string]$SheetName,
        [int]$Row,
        [hashtable]$Values   # column letter -> new value ($null => clear cell)
    )
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($col in $Values.Keys) {
        $cell = $ws.Range("$col$Row")
        $val = $Values[$col]
        if ($null -eq $val) {
            $cell.ClearContents()
        } else {
            $cell.Value = $val
        }
    }
}

# ---------------- ALC ----------------
Set-Row "ALC" 76 @{
    H = 43481504; I = 47622320; J = 2900; K = 47622320; L = 2900
    M = -47622005; N = -3530
}
Set-Row "ALC" 79 @{
    H = 43481504; I = 47622320; J = 2900; K = 47622320; L = 2900
    M = -47621228; N = -5084
}
Set-Row "ALC" 135 @{
    H = 15152226; I = 460.2; J = 62501496; K = 4141.8; L = 562513464
    M = -1606.8; N = -562518534
}

# ---------------- ARM ----------------
Set-Row "ARM" 2 @{
    H = 3136.476; I = 792.3333; J = 11731.667; K = 792.3333; L = 11731.667
    M = -679.3333; N = -11957.667
}
Set-Row "ARM" 32 @{
    H = 5804.8706; I = 4128.568; K = 4128.568; M = -3841.568
}
Set-Row "ARM" 39 @{
    H = 30894.889; I = 4525; J = 51990.8; K = 4525; L = 51990.8
    M = -4005; N = -53030.8
}
Set-Row "ARM" 63 @{
    H = 2438; I = 2438; K = 2438; M = -1752
}
Set-Row "ARM" 66 @{
    H = 2438; I = 2438; K = 12190; M = -8758
}
Set-Row "ARM" 74 @{
    H = 42532.168; I = 59552.47; J = 1197.1428; K = 59552.47; L = 1197.1428
    M = -58678.47; N = -2945.1428
}
Set-Row "ARM" 77 @{
    H = 42532.168; I = 59552.47; J = 1197.1428; K = 297762.35; L = 5985.714
    M = -293394.35; N = -14721.714
}
Set-Row "ARM" 116 @{
    H = 3136.476; I = 792.3333; J = 11731.667; K = 792.3333; L = 11731.667
    M = 1501.6667; N = -16319.667
}

# ---------------- BSM ----------------
Set-Row "BSM" 3 @{
    H = 3136.476; I = 792.3333; J = 11731.667; K = 792.3333; L = 11731.667
    M = -678.3333; N = -11959.667
}
Set-Row "BSM" 7 @{
    H = 1250; I = 0; J = 1250; K = 0; L = 1250
    M = $null; N = -1476
}
Set-Row "BSM" 105 @{
    H = 797530; I = 1328005; J = 1817.5; K = 1328005; L = 1817.5
    M = -1326258; N = -5311.5
}

# ---------------- CRP ----------------
Set-Row "CRP" 35 @{
    H = 26002.9; I = 789; J = 51216.8; K = 789; L = 51216.8
    M = -495; N = -51804.8
}
Set-Row "CRP" 132 @{
    H = 49119.523; I = 101137.4; J = 1830.5454; K = 303412.2; L = 5491.6362
    M = -300882.2; N = -10551.6362
}
Set-Row "CRP" 141 @{
    H = 39475.23; J = 39475.23; L = 39475.23; N = -49835.23
}

# ---------------- CUL ----------------
Set-Row "CUL" 11 @{
    H = 588530.4; I = 909260.6; J = 524.8333; K = 2727781.8; L = 1574.4999
    M = -2727641.8; N = -1854.4999
}
Set-Row "CUL" 40 @{
    H = 218.35294; I = 189.33333; J = 288; K = 757.33332; L = 1152
    M = -688.33332; N = -1290
}
Set-Row "CUL" 46 @{
    H = 2861.261; I = 464.66666; J = 3707.1177; K = 1393.99998; L = 11121.3531
    M = -1302.99998; N = -11303.3531
}
Set-Row "CUL" 49 @{
    H = 103; J = 0; L = 0; N = $null
}
Set-Row "CUL" 60 @{
    H = 675.8333; I = 209.8; J = 3006; K = 629.4000000000001; L = 9018
    M = -378.4000000000001; N = -9520
}
Set-Row "CUL" 70 @{
    H = 2113.875; I = 1462.2; J = 3200; K = 4386.6; L = 9600
    M = -4071.6; N = -10230
}
Set-Row "CUL" 73 @{
    H = 2113.875; I = 1462.2; J = 3200; K = 4386.6; L = 9600
    M = -3294.6; N = -11784
}
Set-Row "CUL" 76 @{
    H = 14614.917; I = 4139.8; J = 22097.143; K = 12419.4; L = 66291.429
    M = -12036.4; N = -67057.429
}
Set-Row "CUL" 79 @{
    H = 14614.917; I = 4139.8; J = 22097.143; K = 12419.4; L = 66291.429
    M = -11093.4; N = -68943.429
}
Set-Row "CUL" 122 @{
    H = 610.7241; I = 377.15384; J = 800.5; K = 3394.38456; L = 7204.5
    M = -944.38456; N = -12104.5
}
Set-Row "CUL" 125 @{
    H = 5975; I = 0; J = 5975; K = 0; L = 17925
    M = $null; N = -27765
}

# ---------------- GSM ----------------
Set-Row "GSM" 70 @{
    H = 4171529.2; I = 5887047; J = 5271.4287; K = 5887047; L = 5271.4287
    M = -5886777; N = -5811.4287
}
Set-Row "GSM" 73 @{
    H = 4171529.2; I = 5887047; J = 5271.4287; K = 5887047; L = 5271.4287
    M = -5886111; N = -7143.4287
}
Set-Row "GSM" 140 @{
    H = 52968.117; J = 52968.117; L = 52968.117; N = -63328.117
}

# ---------------- WVR ----------------
Set-Row "WVR" 76 @{
    H = 16125; J = 16125; L = 16125; N = -16755
}
Set-Row "WVR" 79 @{
    H = 16125; J = 16125; L = 16125; N = -18309
}
